# Auto-generated edit script applying the diff to Aegis_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3880
$ws.Range("I76").Value = 3316.6667
$ws.Range("J76").Value = 4725
$ws.Range("K76").Value = 3316.6667
$ws.Range("L76").Value = 4725
$ws.Range("M76").Value = -3001.6667
$ws.Range("N76").Value = -5355

$ws.Range("H79").Value = 3880
$ws.Range("I79").Value = 3316.6667
$ws.Range("J79").Value = 4725
$ws.Range("K79").Value = 3316.6667
$ws.Range("L79").Value = 4725
$ws.Range("M79").Value = -2224.6667
$ws.Range("N79").Value = -6909

$ws.Range("H132").Value = 10425572
$ws.Range("I132").Value = 11373146
$ws.Range("J132").Value = 2250
$ws.Range("K132").Value = 34119438
$ws.Range("L132").Value = 6750
$ws.Range("M132").Value = -34116908
$ws.Range("N132").Value = -11810

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 73219.86
$ws.Range("I2").Value = 1666.4445
$ws.Range("J2").Value = 202016
$ws.Range("K2").Value = 1666.4445
$ws.Range("L2").Value = 202016
$ws.Range("M2").Value = -1553.4445
$ws.Range("N2").Value = -202242

$ws.Range("H32").Value = 23580.621
$ws.Range("I32").Value = 3725.2424
$ws.Range("K32").Value = 3725.2424
$ws.Range("M32").Value = -3438.2424

$ws.Range("H45").Value = 1630.1177
$ws.Range("I45").Value = 1195.6666
$ws.Range("K45").Value = 1195.6666
$ws.Range("M45").Value = -818.6666

$ws.Range("H116").Value = 73219.86
$ws.Range("I116").Value = 1666.4445
$ws.Range("J116").Value = 202016
$ws.Range("K116").Value = 1666.4445
$ws.Range("L116").Value = 202016
$ws.Range("M116").Value = 627.5554999999999
$ws.Range("N116").Value = -206604

$ws.Range("H133").Value = 65000
$ws.Range("J133").Value = 65000
$ws.Range("L133").Value = 65000
$ws.Range("N133").Value = -70060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 73219.86
$ws.Range("I3").Value = 1666.4445
$ws.Range("J3").Value = 202016
$ws.Range("K3").Value = 1666.4445
$ws.Range("L3").Value = 202016
$ws.Range("M3").Value = -1552.4445
$ws.Range("N3").Value = -202244

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 2549.7
$ws.Range("I86").Value = 2000
$ws.Range("K86").Value = 2000
$ws.Range("M86").Value = -877

$ws.Range("H89").Value = 2549.7
$ws.Range("I89").Value = 2000
$ws.Range("K89").Value = 10000
$ws.Range("M89").Value = -4384

$ws.Range("H99").Value = 8344.235000000001
$ws.Range("I99").Value = 2596.5
$ws.Range("J99").Value = 11479.363
$ws.Range("K99").Value = 2596.5
$ws.Range("L99").Value = 11479.363
$ws.Range("M99").Value = -1098.5
$ws.Range("N99").Value = -14475.363

$ws.Range("H126").Value = 8344.235000000001
$ws.Range("I126").Value = 2596.5
$ws.Range("J126").Value = 11479.363
$ws.Range("K126").Value = 7789.5
$ws.Range("L126").Value = 34438.089
$ws.Range("M126").Value = -5319.5
$ws.Range("N126").Value = -39378.089

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 313.875
$ws.Range("I11").Value = 217
$ws.Range("J11").Value = 372
$ws.Range("K11").Value = 651
$ws.Range("L11").Value = 1116
$ws.Range("M11").Value = -511
$ws.Range("N11").Value = -1396

$ws.Range("H34").Value = 2257.0715
$ws.Range("I34").Value = 149.5
$ws.Range("J34").Value = 2608.3333
$ws.Range("K34").Value = 448.5
$ws.Range("L34").Value = 7824.999899999999
$ws.Range("M34").Value = -364.5
$ws.Range("N34").Value = -7992.999899999999

$ws.Range("H39").Value = 21967.666
$ws.Range("J39").Value = 32500
$ws.Range("L39").Value = 97500
$ws.Range("N39").Value = -98088

$ws.Range("H55").Value = 7996.1875
$ws.Range("J55").Value = 8489.267
$ws.Range("L55").Value = 25467.801
$ws.Range("N55").Value = -25821.801

$ws.Range("H62").Value = 5000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 15000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -16372

$ws.Range("H65").Value = 5000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 45000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -51864

$ws.Range("H130").Value = 744.5
$ws.Range("I130").Value = 744.5
$ws.Range("K130").Value = 2233.5
$ws.Range("M130").Value = 2786.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 27128.863
$ws.Range("I70").Value = 44201.52
$ws.Range("J70").Value = 4664.8423
$ws.Range("K70").Value = 44201.52
$ws.Range("L70").Value = 4664.8423
$ws.Range("M70").Value = -43931.52
$ws.Range("N70").Value = -5204.8423

$ws.Range("H73").Value = 27128.863
$ws.Range("I73").Value = 44201.52
$ws.Range("J73").Value = 4664.8423
$ws.Range("K73").Value = 44201.52
$ws.Range("L73").Value = 4664.8423
$ws.Range("M73").Value = -43265.52
$ws.Range("N73").Value = -6536.8423

$ws.Range("H102").Value = 2045.6842
$ws.Range("I102").Value = 2033.4333
$ws.Range("J102").Value = 2091.625
$ws.Range("K102").Value = 2033.4333
$ws.Range("L102").Value = 2091.625
$ws.Range("M102").Value = -411.4332999999999
$ws.Range("N102").Value = -5335.625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4627.1763
$ws.Range("I132").Value = 4831.846
$ws.Range("J132").Value = 3962
$ws.Range("K132").Value = 14495.538
$ws.Range("L132").Value = 11886
$ws.Range("M132").Value = -11965.538
$ws.Range("N132").Value = -16946

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 5491.1665
$ws.Range("J54").Value = 5175.4
$ws.Range("L54").Value = 5175.4
$ws.Range("N54").Value = -6215.4

$ws.Range("H113").Value = 513.1429000000001
$ws.Range("I113").Value = 405.7857
$ws.Range("J113").Value = 727.8570999999999
$ws.Range("K113").Value = 1217.3571
$ws.Range("L113").Value = 2183.5713
$ws.Range("M113").Value = 952.6428999999998
$ws.Range("N113").Value = -6523.5713

$ws.Range("H126").Value = 1582.6666
$ws.Range("I126").Value = 1572.84
$ws.Range("J126").Value = 1631.8
$ws.Range("K126").Value = 4718.52
$ws.Range("L126").Value = 4895.4
$ws.Range("M126").Value = -2248.52
$ws.Range("N126").Value = -9835.4

